$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.951.88'
$ws.Range('E2').Value = '  +2.74%  '

$ws.Range('D3').Value = '3.514.17'
$ws.Range('E3').Value = '  +1.73%  '

$ws.Range('E4').Value = '  +0.19%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '603.14'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +3.78%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.14'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +4.14%  '

$ws.Range('E7').Value = '  +0.99%  '

$ws.Range('D8').Value = '3.507.14'
$ws.Range('E8').Value = '  +1.91%  '

$ws.Range('E9').Value = '  -0.03%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.191'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +3.48%  '

$ws.Range('E11').Value = '  +6.64%  '

$ws.Range('E12').Value = '  +2.22%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '45.88'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.13%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000274'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.64%  '

$ws.Range('D15').Value = '4.083.28'
$ws.Range('E15').Value = '  +2.37%  '

$ws.Range('E16').Value = '  -0.17%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '603.92'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.24%  '

$ws.Range('D18').Value = '3.534.32'
$ws.Range('E18').Value = '  +2.51%  '

$ws.Range('D19').Value = '69.944.96'
$ws.Range('E19').Value = '  +3.02%  '

$ws.Range('E20').Value = '  +1.03%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.17'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.47%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.868'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.12%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.21'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -15.45%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '15.53'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.58%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '95.70'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.68%  '

$ws.Range('E26').Value = '  -0.22%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.12%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.58'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.64%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.84'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +5.59%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.97'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.03%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '719.55'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +23.63%  '

$ws.Range('E32').Value = '  -0.32%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '8.08'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.98%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.92'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.65%  '

$ws.Range('E35').Value = '  +0.01%  '

$ws.Range('E36').Value = '  -0.47%  '

$ws.Range('E37').Value = '  +5.40%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '10.65'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.45%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0472'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +9.83%  '

$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '56.63'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.16%  '

$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.31%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.142'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +5.07%  '

$ws.Range('D43').Value = '3.347.34'
$ws.Range('E43').Value = '  -0.12%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.314'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.31%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '32.27'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.01%  '

$ws.Range('B46').Value = 'ThetaToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.89'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +6.22%  '

$ws.Range('B47').Value = 'PEPE'
$ws.Range('C47').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D47').Value = '0.0₃0688'
$ws.Range('E47').Value = '  +1.61%  '

$ws.Range('E48').Value = '  +2.81%  '

$ws.Range('E49').Value = '  +1.50%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.88'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.60%  '

$ws.Range('E51').Value = '  -0.03%  '
